# Updates cryptos list values (price/volume) and reorders two rows
# per upstream data refresh, matching commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.876.28"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.304.69"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "2.302.57"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "59.847.76"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "2.714.49"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "2.306.28"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D32").Value = "0.0₃0722"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.95%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "317.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.566"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "0.0₆0224"
$ws.Range("E49").Value = "  +21.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0213"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
